# Added Crdc login backup codes
# Replace the old backup codes (rows 2-4) with the next batch of codes that
# used to live at rows 17-18, then clear the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the stale codes in rows 2-4 with the two fresh codes that used
# to sit further down the sheet (rows 17-18).
$ws.Range("A2").Value = "8AW7QA18SBTA"
$ws.Range("A3").Value = "EA5XZ049QR7S"
$ws.Range("A4").Value = "PTV3TSFPBF6W"

# Move the remaining two codes up to rows 11-12.
$ws.Range("A11").Value = "NHYK5008HQDA"
$ws.Range("A12").Value = "ZADGNDVPP03M"

# Clear out the old rows 17-21, which are no longer used.
$ws.Range("A13:A21").ClearContents()

$ws.Range("C17").Select()
